$d = $word.ActiveDocument

# Locate the two trailing empty paragraphs at the end of the document body.
# The new content must be inserted between them.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

$insertionPoint = $lastPara.Range
$insertionPoint.Collapse(1)  # wdCollapseStart

$newContentXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Common Vulnerabilities and Exposures (CVE) and Security Databases</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Vulnerabilities vs. Exposures</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Vulnerability</w:t></w:r><w:r><w:t xml:space="preserve">: A </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>weakness</w:t></w:r><w:r><w:t xml:space="preserve"> in a system that can be exploited by a threat.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Exposure</w:t></w:r><w:r><w:t xml:space="preserve">: A </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>mistake</w:t></w:r><w:r><w:t xml:space="preserve"> that creates an opportunity for a threat.</w:t></w:r></w:p><w:p><w:r><w:t>Example:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">A </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>document</w:t></w:r><w:r><w:t xml:space="preserve"> left near an open window is exposed to being </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>blown away</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>The CVE List (Common Vulnerabilities and Exposures)</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">A </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>public dictionary</w:t></w:r><w:r><w:t xml:space="preserve"> of known security flaws.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Created by </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>MITRE Corporation</w:t></w:r><w:r><w:t xml:space="preserve"> in 1999.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Used by organizations to </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>identify and mitigate security risks</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Anyone can report</w:t></w:r><w:r><w:t xml:space="preserve"> a vulnerability, but it must pass a strict </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>review process</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>CVE Review Process &amp; Criteria</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Before assigning a </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>CVE ID</w:t></w:r><w:r><w:t>, vulnerabilities must meet four conditions:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Independence</w:t></w:r><w:r><w:t xml:space="preserve"> – Can be fixed without addressing other flaws.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Recognized Risk</w:t></w:r><w:r><w:t xml:space="preserve"> – Reporter must acknowledge its security impact.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Supporting Evidence</w:t></w:r><w:r><w:t xml:space="preserve"> – Requires documentation or proof.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Single Codebase</w:t></w:r><w:r><w:t xml:space="preserve"> – Affects only one system/version.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>NIST National Vulnerability Database (NVD) &amp; CVSS Scoring</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">NVD </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>analyzes CVEs</w:t></w:r><w:r><w:t xml:space="preserve"> further and assigns a </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>severity score</w:t></w:r><w:r><w:t xml:space="preserve"> using </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>CVSS (Common Vulnerability Scoring System)</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Scoring Scale (0-10):</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>0-3.9</w:t></w:r><w:r><w:t xml:space="preserve"> → Low Risk (Not urgent)</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>4.0-6.9</w:t></w:r><w:r><w:t xml:space="preserve"> → Medium Risk</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>7.0-8.9</w:t></w:r><w:r><w:t xml:space="preserve"> → High Risk</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="1"/><w:numId w:val="8"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t>9.0-10</w:t></w:r><w:r><w:t xml:space="preserve"> → Critical Risk (Immediate attention required)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Why Security Teams Use CVE Lists &amp; CVSS Scores</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:t>Helps</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>prioritize security patches</w:t></w:r><w:r><w:t xml:space="preserve"> and updates.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Provides </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>global insights</w:t></w:r><w:r><w:t xml:space="preserve"> into cybersecurity threats.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Supports organizations in </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>making informed security decisions</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($newContentXml)

Write-Output "Inserted new content. Paragraph count now: $($d.Paragraphs.Count)"
